# Weekly refresh of the "Vega Monumental Concepción - Perejil" price series.
# A new reporting date is inserted at the top of the data block (after the
# header row), pushing the existing rows down by one observation (two rows,
# since each date has a "Primera" and "Segunda" quality row), so the oldest
# observation that falls off the bottom of the table is appended at the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new (blank) rows right before the current first data pair
# (row 82), shifting rows 82:151 down to 84:153.
$ws.Rows("82:83").Insert()

# Seed the new top rows with a copy of the row pair now sitting just below
# them (the previous newest observation), preserving all formatting/values.
$ws.Range("A84:R85").Copy($ws.Range("A82"))

# Only the date changes for the new observation.
$ws.Range("D82").Value = 44719
$ws.Range("D83").Value = 44719
